# "Add files via upload" - refresh the COVID-19 Valais daily figures sheet
# with newly reported rows and move the active selection down from the
# merged title row (A1:M1) to the header row (A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Update the two most recent existing rows with corrected "new cases" counts ---
$ws.Range("C589").Value = 51
$ws.Range("C590").Value = 56

# --- Fill in the next four rows that were previously blank placeholders ---
# Columns L and M are formatted as Text ("@"); assigning a numeric .Value
# straight through a Text-formatted cell stores the digits as text, so we
# flip those cells to General just long enough to write a real number, then
# restore the original Text format (matching every other row's layout).
$lmRange = $ws.Range("L591:M594")
$lmRange.NumberFormat = "General"

# Row 591 (2021-10-08)
$ws.Range("C591").Value = 36
$ws.Range("E591").Value = 2
$ws.Range("F591").Value = 1
$ws.Range("G591").Value = 9
$ws.Range("L591").Value = 0
$ws.Range("M591").Value = 0

# Row 592 (2021-10-09)
$ws.Range("C592").Value = 25
$ws.Range("E592").Value = 2
$ws.Range("F592").Value = 1
$ws.Range("G592").Value = 7
$ws.Range("L592").Value = 0
$ws.Range("M592").Value = 0

# Row 593 (2021-10-10)
$ws.Range("C593").Value = 22
$ws.Range("E593").Value = 3
$ws.Range("F593").Value = 1
$ws.Range("G593").Value = 7
$ws.Range("L593").Value = 0
$ws.Range("M593").Value = 0

# Row 594 (2021-10-11)
$ws.Range("C594").Value = 1
$ws.Range("E594").Value = 3
$ws.Range("F594").Value = 1
$ws.Range("G594").Value = 7
$ws.Range("L594").Value = 0
$ws.Range("M594").Value = 0

$lmRange.NumberFormat = "@"

# --- Move the selection off the merged title banner onto the header row ---
$ws.Range("A2").Select()
